# Update the shipment/receiving log on Sheet1 with the new batch of rows
# (old MBK / 01-10-68 invoice pair 6500118096-097 is replaced by the
# BKK/FPR pair 6500128102-103 dated 01/11/25-11/11/25).
#
# Every data cell in this sheet is formatted as Text (numFmtId 49), so
# string assignments below are intentional even for the numeric-looking
# Qty/UnitCost/Date values - only the Invoice column (D) is a real number.
#
# The order in which the cells are written mirrors the order the values
# were entered into the sheet (Dept/Qty/UnitCost per row, then the Date
# column, then the Invoice numbers), so the shared-string table is
# rebuilt in that same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BKK
$ws.Range("A2").Value = "BKK"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "5,014.85"

# Row 3 - FPR
$ws.Range("A3").Value = "FPR"
$ws.Range("G3").Value = "935.97"

# Dates for both rows
$ws.Range("B2").Value = "011125"
$ws.Range("B3").Value = "111125"

# Invoice numbers (real numbers, not text) and the repeated Qty text
$ws.Range("D2").Value = 6500128102
$ws.Range("D3").Value = 6500128103
$ws.Range("F3").Value = "1"

# Column D widened from 11 to 12 characters
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666

# Leave the selection where the user left it when they saved
$null = $ws.Range("E9").Select()
